# Partial build of a board
# Updated inventory to reflect components used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Stock count for 8M resonator (X101) went from 6 to 5
$ws.Range("B12").Value = 5

# Row 21: Stock count for microB USB Female connector went from 71 to 70
$ws.Range("B21").Value = 70

# Row 36: Stock count for ATMEGA328P-AU microcontroller went from 5 to 4
$ws.Range("B36").Value = 4

# Row 40: C105 capacitor value corrected from "10n" to "0.01u" (same value,
# consistent naming with other 0.01u capacitors), and stock count updated
# from 9 to 14
$ws.Range("C40").Value = "0.01u"
$ws.Range("B40").Value = 14

# Row 42: Stock count for LP2985-33DBVR regulator went from 4 to 3
$ws.Range("B42").Value = 3

# Update the view state to reflect where the user was working
$ws.Range("B22").Select()

$wb.Save()
